# cs-en-us-006pct.xlsx weekly refresh: new Volume/date header + updated crime-stat rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text: Volume/Number and "Report Covering the Week" date range ---
# A8 shared string = "Volume 32   Number  46" -> replace the "46" run with "47"
$ws.Range("A8").Characters(21, 2).Text = "47"
# C9 shared string = "Report Covering the Week  11/10/2025  Through  11/16/2025"
$ws.Range("C9").Characters(27, 10).Text = "11/17/2025"
$ws.Range("C9").Characters(48, 10).Text = "11/23/2025"

# --- Cells whose underlying style flips between the blank-placeholder style (13)
#     and a numeric style (14/15): copy a same-style donor cell first so the xf/
#     number-format matches Excel’s own "paste like neighbour" behaviour, then set
#     the real value. Donors are row 33, a row untouched by this edit. ---
$ws.Range("F33").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K33").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("F33").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 7
$ws.Range("K33").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -71.428571428571
$ws.Range("C33").Copy($ws.Range("C20"))
$ws.Range("F33").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 1
$ws.Range("K33").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -100
$ws.Range("F33").Copy($ws.Range("G20"))
$ws.Range("G20").Value = 1
$ws.Range("K33").Copy($ws.Range("H20"))
$ws.Range("H20").Value = 400
$ws.Range("C33").Copy($ws.Range("C22"))
$ws.Range("F33").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K33").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100
$ws.Range("C33").Copy($ws.Range("C28"))
$ws.Range("F33").Copy($ws.Range("D31"))
$ws.Range("D31").Value = 1
$ws.Range("K33").Copy($ws.Range("E31"))
$ws.Range("E31").Value = -100
$ws.Range("C33").Copy($ws.Range("F31"))
$ws.Range("F33").Copy($ws.Range("G31"))
$ws.Range("G31").Value = 1
$ws.Range("K33").Copy($ws.Range("H31"))
$ws.Range("H31").Value = -100

# --- Remaining cells: same style before/after, just refresh the numbers ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 25
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 115
$ws.Range("J16").Value = 151
$ws.Range("K16").Value = -23.841059602649
$ws.Range("L16").Value = -34.659090909090
$ws.Range("M16").Value = -13.533834586466
$ws.Range("N16").Value = -85.237483953786
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -18.181818181818
$ws.Range("I17").Value = 142
$ws.Range("J17").Value = 158
$ws.Range("K17").Value = -10.126582278481
$ws.Range("L17").Value = -21.546961325966
$ws.Range("M17").Value = 65.116279069767
$ws.Range("N17").Value = -49.823321554770
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 151
$ws.Range("J18").Value = 226
$ws.Range("K18").Value = -33.185840707964
$ws.Range("L18").Value = -43.233082706766
$ws.Range("M18").Value = -6.790123456790
$ws.Range("N18").Value = -79.539295392953
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 22
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 109
$ws.Range("G19").Value = 129
$ws.Range("H19").Value = -15.503875968992
$ws.Range("I19").Value = 957
$ws.Range("J19").Value = 1077
$ws.Range("K19").Value = -11.142061281337
$ws.Range("L19").Value = -20.646766169154
$ws.Range("M19").Value = 1.162790697674
$ws.Range("N19").Value = -58.927038626609
$ws.Range("J20").Value = 31
$ws.Range("K20").Value = -16.129032258064
$ws.Range("L20").Value = -40.909090909090
$ws.Range("M20").Value = -36.585365853658
$ws.Range("N20").Value = -96.275071633237
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -24.242424242424
$ws.Range("F21").Value = 139
$ws.Range("G21").Value = 165
$ws.Range("H21").Value = -15.757575757575
$ws.Range("I21").Value = 1402
$ws.Range("J21").Value = 1651
$ws.Range("K21").Value = -15.081768625075
$ws.Range("L21").Value = -25.504782146652
$ws.Range("M21").Value = 1.594202898550
$ws.Range("N21").Value = -71.033057851239
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 37
$ws.Range("K22").Value = 8.108108108108
$ws.Range("M22").Value = -11.111111111111
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = -59.183673469387
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 179
$ws.Range("H24").Value = -42.458100558659
$ws.Range("I24").Value = 1385
$ws.Range("J24").Value = 1700
$ws.Range("K24").Value = -18.529411764705
$ws.Range("L24").Value = -26.172707889125
$ws.Range("M24").Value = 2.592592592592
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 38
$ws.Range("E25").Value = -65.789473684210
$ws.Range("F25").Value = 63
$ws.Range("G25").Value = 140
$ws.Range("H25").Value = -55
$ws.Range("I25").Value = 1041
$ws.Range("J25").Value = 1334
$ws.Range("K25").Value = -21.964017991004
$ws.Range("L25").Value = -27.808599167822
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 31
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -18.421052631578
$ws.Range("I26").Value = 337
$ws.Range("J26").Value = 353
$ws.Range("K26").Value = -4.532577903682
$ws.Range("L26").Value = -13.367609254498
$ws.Range("M26").Value = 43.404255319148
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = 0
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 12
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 69
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = 6.153846153846
$ws.Range("L28").Value = 18.965517241379
$ws.Range("N29").Value = -63.636363636363
$ws.Range("N30").Value = -72.727272727272
$ws.Range("J31").Value = 21
$ws.Range("K31").Value = -47.619047619047

Write-Output "edits applied"
